# Remove the "SLB-Private" classification watermark text boxes from every
# section footer. Each footer's single paragraph currently holds one run
# that wraps an mc:AlternateContent drawing (the classification shape);
# deleting the Shape leaves the footer paragraph empty (just the Footer
# style), matching the target edit.

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    # Footers collection: 1 = primary, 2 = first page, 3 = even page.
    for ($i = 1; $i -le 3; $i++) {
        $ftr = $sec.Footers.Item($i)
        while ($ftr.Shapes.Count -gt 0) {
            $ftr.Shapes.Item(1).Delete()
        }
    }
}
